$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new group-member rows under the existing header row.
$ws.Range("A2").Value = 22599
$ws.Range("B2").Value = 150220078
$ws.Range("C2").Value = "Alperen Akbaş"

$ws.Range("A3").Value = 22599
$ws.Range("B3").Value = 150230046
$ws.Range("C3").Value = "Ömer Faruk Ekmekçi"

# Match the printer/page setup added in the commit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Reflect the new active cell / selection left behind by the edit.
[void]$ws.Range("C3").Select()
